# Update "why we sleep" time records sheet with new rows of activity data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: UK visa, 20 minutes
$ws.Range("B36").Value = "UK签证"
$ws.Range("D36").Value = 20

# Row 37: transportation / visa, 70 minutes
$ws.Range("B37").Value = "交通"
$ws.Range("C37").Value = "签证"
$ws.Range("D37").Value = 70

# Row 38: phone call / grandfather, 5 minutes
$ws.Range("B38").Value = "电话"
$ws.Range("C38").Value = "外公"
$ws.Range("D38").Value = 5

# Row 39: eating, 20 minutes
$ws.Range("B39").Value = "吃饭"
$ws.Range("D39").Value = 20

# Row 40: course / teaching, 47 minutes
$ws.Range("B40").Value = "课程"
$ws.Range("C40").Value = "teaching"
$ws.Range("D40").Value = 47

# Row 41 (new): check Birmingham accommodation, 15 minutes
$ws.Range("A41").Value = 44949
$ws.Range("B41").Value = "查伯明翰住宿"
$ws.Range("D41").Value = 15

# Row 42 (new): nap, 40 minutes
$ws.Range("A42").Value = 44949
$ws.Range("B42").Value = "午睡"
$ws.Range("D42").Value = 40

# Row 43 (new): reading
$ws.Range("A43").Value = 44949
$ws.Range("B43").Value = "阅读"

# Row 44 (new): just the date
$ws.Range("A44").Value = 44949

# Row 45 (new): just the date
$ws.Range("A45").Value = 44949

# Row 46 (new): just the date
$ws.Range("A46").Value = 44949

# Update current selection to match the author's last edit position
$ws.Range("C43").Select()
